# Weekly data update: a new observation row is inserted at row 49 of the
# "Hortaliza, Macroferia Regional de Talca - Cilantro" sheet, pushing the
# previously existing rows 49-131 down to 50-132 (dimension grows from
# A1:R131 to A1:R132). The new row carries a fresh week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 49; everything below shifts down one row,
# carrying its original values and formatting with it.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly record.
$ws.Cells.Item(49, 1).Value = 5
$ws.Cells.Item(49, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(49, 3).Value = "Maule"
$ws.Cells.Item(49, 4).Value = 45203
$ws.Cells.Item(49, 5).Value = 7
$ws.Cells.Item(49, 6).Value = 100112040
$ws.Cells.Item(49, 7).Value = "Cilantro"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 150
$ws.Cells.Item(49, 11).Value = 9000
$ws.Cells.Item(49, 12).Value = 9000
$ws.Cells.Item(49, 13).Value = 9000
$ws.Cells.Item(49, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 250
$ws.Cells.Item(49, 17).Value = 36
$ws.Cells.Item(49, 18).Value = "Hortaliza"
